$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.183899760246277
$ws.Range("B1").Value = 2.174247741699219
$ws.Range("C1").Value = 3.413800954818726
$ws.Range("D1").Value = 1.871991157531738
$ws.Range("E1").Value = 1.013299226760864
